$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text (e.g. "65.05") would otherwise be
# auto-converted to a number by Excel; force them to remain plain text so
# they match the original inline-string cell type/formatting.
$textCells = @("D5", "D8", "D10", "D11", "D16", "D19", "D20", "D22", "D24", "D30", "D37", "D42", "D43", "D45", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.723.40"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.599.30"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "211.62"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.0619"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").Value = "19.62"
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").Value = "1.822.94"
$ws.Range("D13").Value = "1.594.65"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "65.05"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("E17").Value = "  -3.21%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "208.77"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "7.15"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").Value = "2.23"
$ws.Range("E22").Value = "  -3.73%  "
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").Value = "143.95"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("D30").Value = "1.15"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").Value = "1.275.85"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("E34").Value = "  +16.61%  "
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("D37").Value = "0.589"
$ws.Range("E37").Value = "  -4.36%  "
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").Value = "0.777"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").Value = "62.60"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "1.734.82"
$ws.Range("D45").Value = "90.44"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("D49").Value = "7.56"
$ws.Range("E49").Value = "  +3.77%  "
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("E51").Value = "  +1.56%  "
